$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly data rows (2,3,4) got rotated: new row2 = old row4,
# new row3 = old row2, new row4 = old row3 for the columns that
# vary week to week (D, M, N, O, P, R, S).

$oldD2 = $ws.Range("D2").Value2
$oldD3 = $ws.Range("D3").Value2
$oldD4 = $ws.Range("D4").Value2

$oldM2 = $ws.Range("M2").Value2
$oldM3 = $ws.Range("M3").Value2
$oldM4 = $ws.Range("M4").Value2

$oldN2 = $ws.Range("N2").Value2
$oldN3 = $ws.Range("N3").Value2
$oldN4 = $ws.Range("N4").Value2

$oldO2 = $ws.Range("O2").Value2
$oldO3 = $ws.Range("O3").Value2
$oldO4 = $ws.Range("O4").Value2

$oldP2 = $ws.Range("P2").Value2
$oldP3 = $ws.Range("P3").Value2
$oldP4 = $ws.Range("P4").Value2

$oldR2 = $ws.Range("R2").Value2
$oldR3 = $ws.Range("R3").Value2
$oldR4 = $ws.Range("R4").Value2

$oldS2 = $ws.Range("S2").Value2
$oldS3 = $ws.Range("S3").Value2
$oldS4 = $ws.Range("S4").Value2

# Row 2 <- old row 4
$ws.Range("D2").Value = $oldD4
$ws.Range("M2").Value = $oldM4
$ws.Range("N2").Value = $oldN4
$ws.Range("O2").Value = $oldO4
$ws.Range("P2").Value = $oldP4
$ws.Range("R2").Value = $oldR4
$ws.Range("S2").Value = $oldS4

# Row 3 <- old row 2
$ws.Range("D3").Value = $oldD2
$ws.Range("M3").Value = $oldM2
$ws.Range("N3").Value = $oldN2
$ws.Range("O3").Value = $oldO2
$ws.Range("P3").Value = $oldP2
$ws.Range("R3").Value = $oldR2
$ws.Range("S3").Value = $oldS2

# Row 4 <- old row 3
$ws.Range("D4").Value = $oldD3
$ws.Range("M4").Value = $oldM3
$ws.Range("N4").Value = $oldN3
$ws.Range("O4").Value = $oldO3
$ws.Range("P4").Value = $oldP3
$ws.Range("R4").Value = $oldR3
$ws.Range("S4").Value = $oldS3

$wb.Save()
